$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add weight to household surplus and fix plotting functions:
# Updated evaluation metrics (FiT, Price, Unmet Demand, Wasted Surplus,
# Household Surplus) recomputed for each RE target row.

$ws.Range("B2").Value = [double]"0.02"
$ws.Range("C2").Value = [double]"0.31"
$ws.Range("D2").Value = [double]"0.1456984653693756"
$ws.Range("E2").Value = [double]"0"
$ws.Range("F2").Value = [double]"1494365305.884841"

$ws.Range("B3").Value = [double]"0.02"
$ws.Range("C3").Value = [double]"0.31"
$ws.Range("D3").Value = [double]"0.1448620987924807"
$ws.Range("E3").Value = [double]"2.625002048089096E-17"
$ws.Range("F3").Value = [double]"1494370764.077965"

$ws.Range("B4").Value = [double]"0.02"
$ws.Range("C4").Value = [double]"0.31"
$ws.Range("D4").Value = [double]"0.1456984653691645"
$ws.Range("E4").Value = [double]"0"
$ws.Range("F4").Value = [double]"1494365305.884843"

$ws.Range("B5").Value = [double]"0.02"
$ws.Range("C5").Value = [double]"0.31"
$ws.Range("D5").Value = [double]"0.1456984653690596"
$ws.Range("E5").Value = [double]"0"
$ws.Range("F5").Value = [double]"1494365305.884843"

$ws.Range("B6").Value = [double]"0.02"
$ws.Range("C6").Value = [double]"0.31"
$ws.Range("D6").Value = [double]"0.1448620987924807"
$ws.Range("E6").Value = [double]"2.625002048089096E-17"
$ws.Range("F6").Value = [double]"1494370764.077965"

$ws.Range("B7").Value = [double]"0.02"
$ws.Range("C7").Value = [double]"0.31"
$ws.Range("D7").Value = [double]"0.1448620987924808"
$ws.Range("E7").Value = [double]"1.852942622180538E-17"
$ws.Range("F7").Value = [double]"1494370764.077965"

$ws.Range("B8").Value = [double]"0"
$ws.Range("C8").Value = [double]"0.31"
$ws.Range("D8").Value = [double]"0.2269628541449643"
$ws.Range("E8").Value = [double]"1.544118851817115E-17"
$ws.Range("F8").Value = [double]"1493530195.173004"
